$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 3.607991271660467
$ws.Range("C2").Value = 1.09648962797047
$ws.Range("D2").Value = 0.07794152141862298
$ws.Range("E2").Value = 0.4058520619764607
$ws.Range("G2").Value = 0.002483511391034103
$ws.Range("N2").Value = 2.622068973487472
$ws.Range("B3").Value = 3.283876603269334
$ws.Range("C3").Value = 0.9801505167470168
$ws.Range("D3").Value = 0.07070513608638862
$ws.Range("E3").Value = 0.3531709394833342
$ws.Range("G3").Value = 0.002494053275355787
$ws.Range("N3").Value = 2.529050544663647
$ws.Range("B4").Value = 3.087689595922882
$ws.Range("C4").Value = 0.9095208016346419
$ws.Range("D4").Value = 0.06632147501564134
$ws.Range("E4").Value = 0.3210618987568523
$ws.Range("G4").Value = 0.002500829857046504
$ws.Range("N4").Value = 2.472358693095117
$ws.Range("B5").Value = 3.008422821027807
$ws.Range("C5").Value = 0.8809298984854763
$ws.Range("D5").Value = 0.06454944845705768
$ws.Range("E5").Value = 0.308031061923657
$ws.Range("G5").Value = 0.002503668223982367
$ws.Range("N5").Value = 2.449355476142642
$ws.Range("B6").Value = 2.995300953634569
$ws.Range("C6").Value = 0.8761936480175905
$ws.Range("D6").Value = 0.06425605394566958
$ws.Range("E6").Value = 0.305870386617741
$ws.Range("G6").Value = 0.002504144187872324
$ws.Range("N6").Value = 2.445541599921569
$ws.Range("B7").Value = 3.086617855846441
$ws.Range("C7").Value = 0.9091344540868249
$ws.Range("D7").Value = 0.06629751953786922
$ws.Range("E7").Value = 0.3208859498670762
$ws.Range("G7").Value = 0.002500867824367341
$ws.Range("N7").Value = 2.472048071152273
$ws.Range("B8").Value = 3.495633757773817
$ws.Range("C8").Value = 1.056202728160429
$ws.Range("D8").Value = 0.07543368867750644
$ws.Range("E8").Value = 0.3876347712616734
$ws.Range("G8").Value = 0.002487083451478061
$ws.Range("N8").Value = 2.589904448452586
$ws.Range("B9").Value = 4.321441842643708
$ws.Range("C9").Value = 1.351479903768166
$ws.Range("D9").Value = 0.09385111391667067
$ws.Range("E9").Value = 0.5206879776510505
$ws.Range("G9").Value = 0.002462441957714907
$ws.Range("N9").Value = 2.824694521077987
$ws.Range("B10").Value = 4.944605965607593
$ws.Range("C10").Value = 1.573351480338943
$ws.Range("D10").Value = 0.1077309888957814
$ws.Range("E10").Value = 0.6201586605249929
$ws.Range("G10").Value = 0.002445764728262779
$ws.Range("N10").Value = 2.999910861485063
$ws.Range("B11").Value = 5.232120058627743
$ws.Range("C11").Value = 1.675525565567227
$ws.Range("D11").Value = 0.1141307638638693
$ws.Range("E11").Value = 0.6658742039239627
$ws.Range("G11").Value = 0.002438481088303795
$ws.Range("N11").Value = 3.080318195923383
$ws.Range("B12").Value = 5.341611437104234
$ws.Range("C12").Value = 1.714409006407095
$ws.Range("D12").Value = 0.1165673297405618
$ws.Range("E12").Value = 0.6832600969204918
$ws.Range("G12").Value = 0.002435765994697693
$ws.Range("N12").Value = 3.11087581487601
$ws.Range("B13").Value = 5.318002523836185
$ws.Range("C13").Value = 1.706025987262137
$ws.Range("D13").Value = 0.1160419758960813
$ws.Range("E13").Value = 0.6795123032104726
$ws.Range("G13").Value = 0.00243634883151349
$ws.Range("N13").Value = 3.104289696420665
$ws.Range("B14").Value = 5.24111542300011
$ws.Range("C14").Value = 1.678720595848631
$ws.Range("D14").Value = 0.1143309542926119
$ws.Range("E14").Value = 0.6673030169911414
$ws.Range("G14").Value = 0.002438256855917479
$ws.Range("N14").Value = 3.082829954159308
$ws.Range("B15").Value = 5.194101170653539
$ws.Range("C15").Value = 1.662020711927994
$ws.Range("D15").Value = 0.1132846350843266
$ws.Range("E15").Value = 0.6598343920336731
$ws.Range("G15").Value = 0.002439431168881801
$ws.Range("N15").Value = 3.069699703873425
$ws.Range("B16").Value = 4.92590045725683
$ws.Range("C16").Value = 1.566700330973845
$ws.Range("D16").Value = 0.107314540023566
$ws.Range("E16").Value = 0.6171810441194765
$ws.Range("G16").Value = 0.002446246784535164
$ws.Range("N16").Value = 2.994670859619561
$ws.Range("B17").Value = 4.762426308702743
$ws.Range("C17").Value = 1.508552331651572
$ws.Range("D17").Value = 0.1036745882182402
$ws.Range("E17").Value = 0.5911391009289986
$ws.Range("G17").Value = 0.002450505177481154
$ws.Range("N17").Value = 2.948828397422801
$ws.Range("B18").Value = 4.668777731438581
$ws.Range("C18").Value = 1.475223355774915
$ws.Range("D18").Value = 0.1015890071188039
$ws.Range("E18").Value = 0.5762038495583113
$ws.Range("G18").Value = 0.002452983034488977
$ws.Range("N18").Value = 2.922526448895354
$ws.Range("B19").Value = 4.637133680761281
$ws.Range("C19").Value = 1.463958293886776
$ws.Range("D19").Value = 0.1008842198563116
$ws.Range("E19").Value = 0.5711542599483437
$ws.Range("G19").Value = 0.00245382691151039
$ws.Range("N19").Value = 2.913632046536179
$ws.Range("B20").Value = 4.779789091841963
$ws.Range("C20").Value = 1.514730161349689
$ws.Range("D20").Value = 0.104061231408096
$ws.Range("E20").Value = 0.5939067671176588
$ws.Range("G20").Value = 0.002450048913805592
$ws.Range("N20").Value = 2.953701574882871
$ws.Range("B21").Value = 5.26368201577111
$ws.Range("C21").Value = 1.686735521255514
$ws.Range("D21").Value = 0.114833160660794
$ws.Range("E21").Value = 0.6708871003298356
$ws.Range("G21").Value = 0.002437695257922277
$ws.Range("N21").Value = 3.08913017159108
$ws.Range("B22").Value = 5.583542398614554
$ws.Range("C22").Value = 1.800278657171361
$ws.Range("D22").Value = 0.1219500366514694
$ws.Range("E22").Value = 0.7216355596672486
$ws.Range("G22").Value = 0.002429872195777015
$ws.Range("N22").Value = 3.178280608481998
$ws.Range("B23").Value = 5.412484475794486
$ws.Range("C23").Value = 1.739570770950081
$ws.Range("D23").Value = 0.1181443339891359
$ws.Range("E23").Value = 0.694507606306189
$ws.Range("G23").Value = 0.002434024732548656
$ws.Range("N23").Value = 3.130637950011078
$ws.Range("B24").Value = 4.771938337189681
$ws.Range("C24").Value = 1.511936850706888
$ws.Range("D24").Value = 0.1038864080413759
$ws.Range("E24").Value = 0.5926553917642821
$ws.Range("G24").Value = 0.002450255098090781
$ws.Range("N24").Value = 2.951498245362615
$ws.Range("B25").Value = 4.09529824119727
$ws.Range("C25").Value = 1.270788513246032
$ws.Range("D25").Value = 0.08881068369353784
$ws.Range("E25").Value = 0.4844218904496103
$ws.Range("G25").Value = 0.002468855374965174
$ws.Range("N25").Value = 2.760735756104452
